$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A16").Value = "Bahamas"
$ws.Range("A47").Value = "Congo (Kinshasa)"
$ws.Range("A48").Value = "Congo (Brazzaville)"
$ws.Range("A56").Value = "Czechia"
$ws.Range("A75").Value = "Gambia"
$ws.Range("A193").Value = "Eswatini"
$ws.Range("A197").Value = "Taiwan*"
$ws.Range("A213").Value = "US"
$ws.Range("A220").Value = "West Bank and Gaza"
